# Fix formatting of floating point numbers scraped into column H (Importe),
# converting Spanish/Argentine formatted numbers (e.g. "3.000,00") into
# plain decimal-dot numbers stored as text (e.g. "3000.00").
#
# Also fix a handful of "Razon social" values where a comma was
# mis-scraped; replace the stray commas with periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column H (Importe) values, row 2 through 315, in order ----
$newAmounts = @(
    "3000.00",
    "18000.00",
    "15990.00",
    "3612000.00",
    "1188000.00",
    "85.00",
    "517629.80",
    "1297100.00",
    "43200.00",
    "192993.00",
    "3750.00",
    "620.00",
    "6304.00",
    "10980.00",
    "545980.88",
    "590076.12",
    "20400.00",
    "399.00",
    "23002.20",
    "10944.00",
    "154478.30",
    "4700.00",
    "10671.80",
    "118574.13",
    "26773.00",
    "23600.00",
    "36057.90",
    "132797.00",
    "10300.00",
    "5890.00",
    "2350.00",
    "33000.00",
    "3580.00",
    "2381.00",
    "759.11",
    "288776.00",
    "140741.00",
    "2568.00",
    "36948.71",
    "30995.00",
    "7816.50",
    "919.87",
    "99715.00",
    "1500.00",
    "3054.00",
    "264.48",
    "243.11",
    "351903.00",
    "618610.00",
    "11988.00",
    "650.00",
    "249440.88",
    "12995.98",
    "341867.90",
    "109805.00",
    "2978.63",
    "1476.88",
    "2839.00",
    "28748.00",
    "10486.75",
    "27000.00",
    "113498.09",
    "403240.00",
    "5656.84",
    "150.00",
    "2245.00",
    "7032.50",
    "27635.46",
    "2700.00",
    "24664.00",
    "6227.30",
    "930.00",
    "483.00",
    "83.10",
    "8100.00",
    "707.30",
    "124.56",
    "157063.10",
    "2425.00",
    "2900.00",
    "481140.22",
    "4500.00",
    "396.30",
    "1050.00",
    "2417.09",
    "40.00",
    "46207.34",
    "9800.00",
    "10089.00",
    "55114.22",
    "15460.00",
    "30000.00",
    "5000.00",
    "46800.00",
    "5800.00",
    "1000.00",
    "3567.00",
    "21900.00",
    "44500.00",
    "9800.00",
    "700.00",
    "16000.00",
    "8500.00",
    "7020.00",
    "1226.20",
    "8626.00",
    "1830.00",
    "4680.00",
    "1365.00",
    "19525.00",
    "496.40",
    "53980.00",
    "1540.00",
    "6217.40",
    "1365.00",
    "46785.09",
    "4089.00",
    "3.00",
    "20786.50",
    "42347.51",
    "78.00",
    "1668.00",
    "36382.60",
    "3750.00",
    "180.00",
    "4210.00",
    "5222.77",
    "48796.62",
    "2500.00",
    "74.00",
    "3414.00",
    "420.00",
    "3800.00",
    "71407.50",
    "40223.60",
    "281233.00",
    "3487.12",
    "500.00",
    "6490.00",
    "1060.00",
    "3384.90",
    "19164.00",
    "78193.52",
    "1924.30",
    "38712.50",
    "21229.00",
    "388.60",
    "12000.00",
    "380.82",
    "5639.98",
    "804.00",
    "6800.00",
    "26600.00",
    "7200.00",
    "13900.00",
    "4380.00",
    "10488.00",
    "10371.16",
    "4700.00",
    "8000.00",
    "6000.00",
    "4000.00",
    "6000.00",
    "47520.00",
    "280937.50",
    "2400.00",
    "6000.00",
    "7503.04",
    "5014.14",
    "346.00",
    "5960.00",
    "3695.00",
    "1124.00",
    "38640.00",
    "27495.00",
    "12500.00",
    "8000.00",
    "18000.00",
    "8000.00",
    "6500.00",
    "5000.00",
    "8508.50",
    "6500.00",
    "8000.00",
    "7000.00",
    "5000.00",
    "6000.00",
    "6000.00",
    "5000.00",
    "5000.00",
    "16000.00",
    "6000.00",
    "7000.00",
    "12500.00",
    "6000.00",
    "8000.00",
    "6500.00",
    "1500.00",
    "18000.00",
    "11000.00",
    "8950.00",
    "69607.04",
    "15000.00",
    "6500.00",
    "6000.00",
    "7000.00",
    "92280.90",
    "3500.00",
    "34720.00",
    "6000.00",
    "27480.00",
    "3250.06",
    "13000.00",
    "9075.00",
    "2580.00",
    "4992.00",
    "10170.00",
    "2580.00",
    "51.45",
    "970.00",
    "1392.00",
    "13299.20",
    "7172.44",
    "3154.06",
    "11440.93",
    "132745.00",
    "39640.02",
    "670.70",
    "6020.00",
    "4280.00",
    "53131.00",
    "10781.10",
    "25024.00",
    "3346.24",
    "5500.00",
    "2310.00",
    "6069.75",
    "12460.00",
    "4300.00",
    "276.00",
    "1050.00",
    "2824.11",
    "8570.52",
    "15000.00",
    "30000.00",
    "30000.00",
    "105000.00",
    "30000.00",
    "30000.00",
    "60000.00",
    "120000.00",
    "120000.00",
    "60000.00",
    "60000.00",
    "1694.00",
    "225500.00",
    "9529316.75",
    "10890.00",
    "5315.48",
    "3800.00",
    "7200.00",
    "9750.00",
    "7697914.54",
    "322500.00",
    "14611356.47",
    "190000.00",
    "322500.00",
    "172000.00",
    "322500.00",
    "322500.00",
    "430000.00",
    "545400.00",
    "878400.00",
    "172000.00",
    "284950.00",
    "322500.00",
    "322500.00",
    "172000.00",
    "430000.00",
    "719500.00",
    "580500.00",
    "322500.00",
    "597000.00",
    "322500.00",
    "172000.00",
    "179920.00",
    "322500.00",
    "273525.00",
    "1865000.00",
    "415000.00",
    "231052.00",
    "9350.00",
    "21798803.28",
    "321478.90",
    "2268562.22",
    "27930.00",
    "1391.77",
    "3450.00",
    "21360.00",
    "224000.00",
    "12705.00",
    "310200.00",
    "7000.00",
    "17100.00",
    "21427.78",
    "10500.00",
    "861.00",
    "1722.00",
    "53592.50",
    "395900.00",
    "34600.00",
    "4900.00",
    "39555.00",
    "9550.00"
)

$startRow = 2
$endRow = 315
$rngH = $ws.Range("H" + $startRow + ":H" + $endRow)

# Force the range to Text format first so Excel stores the values as
# literal strings (shared strings) rather than re-parsing them as
# numbers (which would strip the trailing zeros / change the cell type).
$rngH.NumberFormat = "@"

for ($i = 0; $i -lt $newAmounts.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 8).Value = $newAmounts[$i]
}

# Restore the default "Normal" style so no stray number-format style is
# left applied to these cells.
$rngH.Style = "Normal"

# ---- Fix the four "Razon social" (and matching "Nombre Fantasia") values ----
$ws.Range("E99").Value  = "BOFFELLI. MARIA INES"

$ws.Range("E200").Value = "PARPAGNOLI. PEDRO RICARDO"
$ws.Range("F200").Value = "PARPAGNOLI. PEDRO RICARDO"

$ws.Range("E216").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E236").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

$ws.Range("E217").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
